# The original rows 3-6 (a 4-row block of observation records) get
# cyclically rotated: old row 6's data moves up to row 3, and old rows
# 3, 4, 5 each shift down by one row (to rows 4, 5, 6 respectively).
# Only the columns that actually differ between these rows need to be
# rewritten; the remaining columns (C, K, N, S, T, U, V, W, Y, Z, AA, AB,
# AD, AE, AF, AG, AT, AW, AX, AY, ...) are identical across rows 3-6 and
# are left untouched.
#
# NOTE: reading `.Value` as a bare property yields a descriptor object in
# this host, not the cell's actual content - it must be invoked as a
# method, i.e. `.Value()`.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot "before" values for rows 3-6 (read with Value() before any writes).
$a3 = $ws.Range("A3").Value(); $b3 = $ws.Range("B3").Value(); $d3 = $ws.Range("D3").Value()
$e3 = $ws.Range("E3").Value(); $f3 = $ws.Range("F3").Value(); $g3 = $ws.Range("G3").Value()
$h3 = $ws.Range("H3").Value(); $i3 = $ws.Range("I3").Value(); $j3 = $ws.Range("J3").Value()
$p3 = $ws.Range("P3").Value(); $q3 = $ws.Range("Q3").Value(); $r3 = $ws.Range("R3").Value()

$a4 = $ws.Range("A4").Value(); $b4 = $ws.Range("B4").Value(); $d4 = $ws.Range("D4").Value()
$e4 = $ws.Range("E4").Value(); $f4 = $ws.Range("F4").Value(); $g4 = $ws.Range("G4").Value()
$h4 = $ws.Range("H4").Value(); $i4 = $ws.Range("I4").Value(); $j4 = $ws.Range("J4").Value()
$p4 = $ws.Range("P4").Value(); $q4 = $ws.Range("Q4").Value(); $r4 = $ws.Range("R4").Value()

$a5 = $ws.Range("A5").Value(); $b5 = $ws.Range("B5").Value(); $d5 = $ws.Range("D5").Value()
$e5 = $ws.Range("E5").Value(); $f5 = $ws.Range("F5").Value(); $g5 = $ws.Range("G5").Value()
$h5 = $ws.Range("H5").Value(); $i5 = $ws.Range("I5").Value(); $j5 = $ws.Range("J5").Value()
$p5 = $ws.Range("P5").Value(); $q5 = $ws.Range("Q5").Value(); $r5 = $ws.Range("R5").Value()

$a6 = $ws.Range("A6").Value(); $b6 = $ws.Range("B6").Value(); $d6 = $ws.Range("D6").Value()
$e6 = $ws.Range("E6").Value(); $f6 = $ws.Range("F6").Value(); $g6 = $ws.Range("G6").Value()
$h6 = $ws.Range("H6").Value(); $i6 = $ws.Range("I6").Value(); $j6 = $ws.Range("J6").Value()
$p6 = $ws.Range("P6").Value(); $q6 = $ws.Range("Q6").Value(); $r6 = $ws.Range("R6").Value()

# New row 3 <- old row 6
$ws.Range("A3").Value = $a6; $ws.Range("B3").Value = $b6; $ws.Range("D3").Value = $d6
$ws.Range("E3").Value = $e6; $ws.Range("F3").Value = $f6; $ws.Range("G3").Value = $g6
$ws.Range("H3").Value = $h6; $ws.Range("I3").Value = $i6; $ws.Range("J3").Value = $j6
$ws.Range("P3").Value = $p6; $ws.Range("Q3").Value = $q6; $ws.Range("R3").Value = $r6

# New row 4 <- old row 3
$ws.Range("A4").Value = $a3; $ws.Range("B4").Value = $b3; $ws.Range("D4").Value = $d3
$ws.Range("E4").Value = $e3; $ws.Range("F4").Value = $f3; $ws.Range("G4").Value = $g3
$ws.Range("H4").Value = $h3; $ws.Range("I4").Value = $i3; $ws.Range("J4").Value = $j3
$ws.Range("P4").Value = $p3; $ws.Range("Q4").Value = $q3; $ws.Range("R4").Value = $r3

# New row 5 <- old row 4
$ws.Range("A5").Value = $a4; $ws.Range("B5").Value = $b4; $ws.Range("D5").Value = $d4
$ws.Range("E5").Value = $e4; $ws.Range("F5").Value = $f4; $ws.Range("G5").Value = $g4
$ws.Range("H5").Value = $h4; $ws.Range("I5").Value = $i4; $ws.Range("J5").Value = $j4
$ws.Range("P5").Value = $p4; $ws.Range("Q5").Value = $q4; $ws.Range("R5").Value = $r4

# New row 6 <- old row 5
# Old I5 ("3") was stored as literal text, not a number - a bare numeric
# string assigned to a General-formatted cell would be auto-coerced to a
# number, so force it back to text with a leading apostrophe (the normal
# Excel "store as text" convention) to preserve its original text type.
$ws.Range("A6").Value = $a5; $ws.Range("B6").Value = $b5; $ws.Range("D6").Value = $d5
$ws.Range("E6").Value = $e5; $ws.Range("F6").Value = $f5; $ws.Range("G6").Value = $g5
$ws.Range("H6").Value = $h5; $ws.Range("J6").Value = $j5
if ($i5 -eq $null -or $i5 -eq "") {
    $ws.Range("I6").Value = $i5
} else {
    $ws.Range("I6").Value = "'" + $i5
}
$ws.Range("P6").Value = $p5; $ws.Range("Q6").Value = $q5; $ws.Range("R6").Value = $r5
